$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "sheet_4"

# Apply the smaller "comments" font (10pt Calibri, black) to the used range.
$fontRange = $ws.Range("A1:F12")
$fontRange.Font.Size = 10
$fontRange.Font.Name = "Calibri"
$fontRange.Font.Color = 0

# Row 1: comments marker
$ws.Range("A1").Value = "# comments"

# Row 2: headers
$ws.Range("A2").Value = "Header 1"
$ws.Range("B2").Value = "Header 2"
$ws.Range("C2").Value = "Header 3"
$ws.Range("D2").Value = "Header 4"
$ws.Range("E2").Value = "Header 5"
$ws.Range("F2").Value = "Header 6"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "a"
$ws.Range("C3").Value = "test1"
$ws.Range("D3").Value = "my_str1"
$ws.Range("E3").Value = 11
$ws.Range("F3").Value = "test11"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "b"
$ws.Range("C4").Value = "test2"
$ws.Range("D4").Value = "my_str2"
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = "test12"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "c"
$ws.Range("C5").Value = "test3"
$ws.Range("D5").Value = "my_str3"
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = "test13"

# Select the cell that ends up highlighted on this sheet once it becomes active.
$ws.Range("E6").Select()
